$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- "Ativação:" date (row 8, columns B and C) ---------------------------
# The literal text "01/01/2023" looks like a date, so Excel would normally
# auto-convert it to a date serial number. Prefixing with an apostrophe
# forces it to stay literal text (quotePrefix), exactly like typing
# '01/01/2023 into the cell.
$ws.Range("B8").Value = "'01/01/2023"
$ws.Range("C8").Value = "'01/01/2023"

# --- "Critério:" description (row 18, columns B and C) -------------------
$ws.Range("B18").Value = 'Esta é uma disciplina de caráter fundamental, exigindo dedicação individual para assimilação das definições e conceitos. Isto envolve leitura concentrada para fixação dos conceitos teóricos e realização de exercícios numéricos. Duas provas escritas (P1 e P2) serão aplicadas e com pesos iguais. O desenvolvimento do aluno ao longo do curso será aferido e estimulado por meio de discussões sobre um dado tema, porém sem a atribuição de nota, por conta da subjetividade envolvida.'
$ws.Range("C18").Value = 'Esta é uma disciplina de caráter fundamental, exigindo dedicação individual para assimilação das definições e conceitos. Isto envolve leitura concentrada para fixação dos conceitos teóricos e realização de exercícios numéricos. Duas provas escritas (P1 e P2) serão aplicadas e com pesos iguais. O desenvolvimento do aluno ao longo do curso será aferido e estimulado por meio de discussões sobre um dado tema, porém sem a atribuição de nota, por conta da subjetividade envolvida.'

# --- "Norma de recuperação:" text (row 19, columns B and C) --------------
$ws.Range("B19").Value = ': A Nota final (NF) será calculada da seguinte maneira: NF = (0,4*P1 +0,4* P2+ 0,2*NT) / 3'
$ws.Range("C19").Value = ': A Nota final (NF) será calculada da seguinte maneira: NF = (0,4*P1 +0,4* P2+ 0,2*NT) / 3'

# --- "Bibliografia:" text (row 20, columns B and C) -----------------------
$ws.Range("B20").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR) / 2'
$ws.Range("C20").Value = 'A recuperação será feita por meio de uma prova escrita (PR) e a média de recuperação (MR) calculada pela fórmula: MR = (NF + PR) / 2'
